# Updates cryptos list (prices + volume %) per latest scrape.
# Some price cells hold plain decimal text (e.g. "585.44") that Excel would
# otherwise auto-convert to a number, so force those ranges to Text format
# before writing the new value, keeping them as strings like the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textPriceCells = @(
    "D5",
    "D6",
    "D11",
    "D12",
    "D13",
    "D14",
    "D17",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D27",
    "D29",
    "D30",
    "D32",
    "D38",
    "D39",
    "D40",
    "D41",
    "D44",
    "D47",
    "D48",
    "D49",
    "D51"
)
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.848.78"
$ws.Range("E2").Value = "  +4.86%  "
$ws.Range("D3").Value = "3.108.81"
$ws.Range("E3").Value = "  +2.79%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "585.44"
$ws.Range("E5").Value = "  +3.45%  "
$ws.Range("D6").Value = "144.14"
$ws.Range("E6").Value = "  +2.15%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.103.67"
$ws.Range("E8").Value = "  +2.85%  "
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("E10").Value = "  +10.71%  "
$ws.Range("D11").Value = "5.68"
$ws.Range("E11").Value = "  +6.57%  "
$ws.Range("D12").Value = "0.467"
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("D13").Value = "0.0000245"
$ws.Range("E13").Value = "  +5.06%  "
$ws.Range("D14").Value = "35.45"
$ws.Range("E14").Value = "  +3.54%  "
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("D16").Value = "3.622.30"
$ws.Range("E16").Value = "  +2.82%  "
$ws.Range("D17").Value = "7.17"
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("D18").Value = "3.104.81"
$ws.Range("E18").Value = "  +2.82%  "
$ws.Range("D19").Value = "62.776.79"
$ws.Range("E19").Value = "  +4.82%  "
$ws.Range("D20").Value = "463.80"
$ws.Range("E20").Value = "  +5.68%  "
$ws.Range("D21").Value = "14.06"
$ws.Range("E21").Value = "  +2.02%  "
$ws.Range("D22").Value = "0.728"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").Value = "7.52"
$ws.Range("E23").Value = "  +5.08%  "
$ws.Range("D24").Value = "13.36"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("D25").Value = "82.12"
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "2.25"
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("E28").Value = "  +4.45%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "8.28"
$ws.Range("E29").Value = "  +4.93%  "
$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  +8.07%  "
$ws.Range("D32").Value = "26.92"
$ws.Range("E32").Value = "  +2.84%  "
$ws.Range("E33").Value = "  +7.98%  "
$ws.Range("D34").Value = "0.0₃0833"
$ws.Range("E34").Value = "  +4.63%  "
$ws.Range("E35").Value = "  +9.98%  "
$ws.Range("E36").Value = "  +3.08%  "
$ws.Range("E37").Value = "  +1.22%  "
$ws.Range("D38").Value = "3.17"
$ws.Range("E38").Value = "  +12.85%  "
$ws.Range("D39").Value = "50.94"
$ws.Range("E39").Value = "  +3.35%  "
$ws.Range("D40").Value = "8.80"
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("D41").Value = "431.41"
$ws.Range("E41").Value = "  +6.25%  "
$ws.Range("D42").Value = "2.920.34"
$ws.Range("E42").Value = "  +4.67%  "
$ws.Range("E43").Value = "  +3.41%  "
$ws.Range("D44").Value = "0.277"
$ws.Range("E44").Value = "  +8.68%  "
$ws.Range("E45").Value = "  +2.23%  "
$ws.Range("E46").Value = "  +6.10%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").Value = "35.09"
$ws.Range("E47").Value = "  +2.29%  "
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").Value = "0.998"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").Value = "123.41"
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("D51").Value = "24.67"
$ws.Range("E51").Value = "  +4.12%  "
